$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Friday hours for week of 43157 (row 7)
$ws.Range("F7").Value = 6.75

# Add Saturday hours for the same week (row 7), which was previously blank
$ws.Range("G7").Value = 7.75

# Move the active selection to reflect where the user ended up (G9)
$ws.Range("G9").Select()
